# Actualización automática de scrims_actualizado.xlsx (2025-07-25 18:34:04)
# Appends new scrim result rows to the "Open Business" and "Ring of Fire" sheets.

$wb = $excel.ActiveWorkbook

function Add-ScrimRow {
    param(
        $ws,
        [int]$templateRow,
        [int]$targetRow,
        [string[]]$values
    )

    # Carry over the row's look (fills/borders/fonts) from an existing row that
    # already has the right "Equipo 1" / "Equipo 2" formatting, then stamp in
    # the real values for this match.
    $src = $ws.Range("A" + $templateRow + ":N" + $templateRow)
    $dst = $ws.Range("A" + $targetRow + ":N" + $targetRow)
    $src.Copy($dst)

    for ($i = 0; $i -lt $values.Length; $i++) {
        $ws.Cells.Item($targetRow, $i + 1).Value2 = $values[$i]
    }
}

# ---------------------------------------------------------------------------
# Sheet "Open Business": rows 87-92 (dimension grows from A3:N86 to A3:N92)
# ---------------------------------------------------------------------------
$wsOB = $wb.Worksheets.Item("Open Business")

# Template rows already on this sheet using the current (bordered) look.
$ob_T2 = 86   # existing "Equipo 2" row to copy formatting from
$ob_T1 = 43   # "Equipo 1" formatted row borrowed from "Ring of Fire"
$wsRoFForTemplate = $wb.Worksheets.Item("Ring of Fire")

Add-ScrimRow $wsOB $ob_T2 87 @("CROW","DRACO","BARLEY","LUMI","MORTIS","MEG","Equipo 2","IC|Mebius","IC|RamaZR","IC|Nob?","FUT|GeRo","FUT|Nowy297","FUT|MeOw","20250725T162956.000Z")

# Rows 88,89,91 need "Equipo 1" styling, which doesn't exist yet on this sheet,
# so borrow the look from "Ring of Fire" row 43 first, then fix the values.
$srcT1 = $wsRoFForTemplate.Range("A" + $ob_T1 + ":N" + $ob_T1)
$dst88 = $wsOB.Range("A88:N88")
$srcT1.Copy($dst88)
$wsOB.Cells.Item(88,1).Value2 = "CROW"
$wsOB.Cells.Item(88,2).Value2 = "DRACO"
$wsOB.Cells.Item(88,3).Value2 = "BARLEY"
$wsOB.Cells.Item(88,4).Value2 = "LUMI"
$wsOB.Cells.Item(88,5).Value2 = "MORTIS"
$wsOB.Cells.Item(88,6).Value2 = "MEG"
$wsOB.Cells.Item(88,7).Value2 = "Equipo 1"
$wsOB.Cells.Item(88,8).Value2 = "IC|Mebius"
$wsOB.Cells.Item(88,9).Value2 = "IC|RamaZR"
$wsOB.Cells.Item(88,10).Value2 = "IC|Nob?"
$wsOB.Cells.Item(88,11).Value2 = "FUT|GeRo"
$wsOB.Cells.Item(88,12).Value2 = "FUT|Nowy297"
$wsOB.Cells.Item(88,13).Value2 = "FUT|MeOw"
$wsOB.Cells.Item(88,14).Value2 = "20250725T162725.000Z"

$dst89 = $wsOB.Range("A89:N89")
$srcT1.Copy($dst89)
$wsOB.Cells.Item(89,1).Value2 = "KENJI"
$wsOB.Cells.Item(89,2).Value2 = "LUMI"
$wsOB.Cells.Item(89,3).Value2 = "GUS"
$wsOB.Cells.Item(89,4).Value2 = "GRAY"
$wsOB.Cells.Item(89,5).Value2 = "JACKY"
$wsOB.Cells.Item(89,6).Value2 = "BEA"
$wsOB.Cells.Item(89,7).Value2 = "Equipo 1"
$wsOB.Cells.Item(89,8).Value2 = "IC|Nob?"
$wsOB.Cells.Item(89,9).Value2 = "IC|RamaZR"
$wsOB.Cells.Item(89,10).Value2 = "IC|Mebius"
$wsOB.Cells.Item(89,11).Value2 = "FUT|MeOw"
$wsOB.Cells.Item(89,12).Value2 = "FUT|Nowy297"
$wsOB.Cells.Item(89,13).Value2 = "FUT|GeRo"
$wsOB.Cells.Item(89,14).Value2 = "20250725T162033.000Z"

Add-ScrimRow $wsOB $ob_T2 90 @("KENJI","LUMI","GUS","GRAY","JACKY","BEA","Equipo 2","IC|Nob?","IC|RamaZR","IC|Mebius","FUT|MeOw","FUT|Nowy297","FUT|GeRo","20250725T161759.000Z")

$dst91 = $wsOB.Range("A91:N91")
$srcT1.Copy($dst91)
$wsOB.Cells.Item(91,1).Value2 = "KENJI"
$wsOB.Cells.Item(91,2).Value2 = "LUMI"
$wsOB.Cells.Item(91,3).Value2 = "GUS"
$wsOB.Cells.Item(91,4).Value2 = "GRAY"
$wsOB.Cells.Item(91,5).Value2 = "JACKY"
$wsOB.Cells.Item(91,6).Value2 = "BEA"
$wsOB.Cells.Item(91,7).Value2 = "Equipo 1"
$wsOB.Cells.Item(91,8).Value2 = "IC|Nob?"
$wsOB.Cells.Item(91,9).Value2 = "IC|RamaZR"
$wsOB.Cells.Item(91,10).Value2 = "IC|Mebius"
$wsOB.Cells.Item(91,11).Value2 = "FUT|MeOw"
$wsOB.Cells.Item(91,12).Value2 = "FUT|Nowy297"
$wsOB.Cells.Item(91,13).Value2 = "FUT|GeRo"
$wsOB.Cells.Item(91,14).Value2 = "20250725T161548.000Z"

Add-ScrimRow $wsOB $ob_T2 92 @("CROW","DRACO","BARLEY","LUMI","MORTIS","MEG","Equipo 2","IC|Mebius","IC|RamaZR","IC|Nob?","FUT|GeRo","FUT|Nowy297","FUT|MeOw","20250725T163225.000Z")

# ---------------------------------------------------------------------------
# Sheet "Ring of Fire": rows 57-65 (dimension grows from A3:N56 to A3:N65)
# ---------------------------------------------------------------------------
$wsRoF = $wb.Worksheets.Item("Ring of Fire")

$rof_T2 = 56   # existing "Equipo 2" row to copy formatting from
$rof_T1 = 43   # existing "Equipo 1" row to copy formatting from

Add-ScrimRow $wsRoF $rof_T1 57 @("CROW","BERRY","8-BIT","POCO","SQUEAK","MEG","Equipo 1","TH|LeNain","TH|iKaoss","TH|Zhar","NXT|Rup","NXT|Arthur","NXT|amos","20250725T163036.000Z")
Add-ScrimRow $wsRoF $rof_T1 58 @("CROW","BERRY","8-BIT","POCO","SQUEAK","MEG","Equipo 1","TH|LeNain","TH|iKaoss","TH|Zhar","NXT|Rup","NXT|Arthur","NXT|amos","20250725T162852.000Z")
Add-ScrimRow $wsRoF $rof_T1 59 @("LUMI","DOUG","MEG","DRACO","BO","GRAY","Equipo 1","NXT|Rup","NXT|Arthur","NXT|amos","TH|LeNain","TH|iKaoss","TH|Zhar","20250725T162258.000Z")
Add-ScrimRow $wsRoF $rof_T1 60 @("DRACO","BO","GRAY","LUMI","MEG","DOUG","Equipo 1","TH|LeNain","TH|iKaoss","TH|Zhar","NXT|amos","NXT|Rup","NXT|Arthur","20250725T162013.000Z")
Add-ScrimRow $wsRoF $rof_T1 61 @("DRACO","BO","GRAY","LUMI","MEG","DOUG","Equipo 1","TH|LeNain","TH|iKaoss","TH|Zhar","NXT|amos","NXT|Rup","NXT|Arthur","20250725T161805.000Z")
Add-ScrimRow $wsRoF $rof_T2 62 @("KIT","HANK","BYRON","CHARLIE","FINX","DRACO","Equipo 2","LOUD|FireCrow","LOUD|Edinho","LOUD|KaioDog","Bielz","Tilo🍥","GO|Yichy❦","20250725T162724.000Z")
Add-ScrimRow $wsRoF $rof_T2 63 @("KIT","HANK","BYRON","CHARLIE","FINX","DRACO","Equipo 2","LOUD|FireCrow","LOUD|Edinho","LOUD|KaioDog","Bielz","Tilo🍥","GO|Yichy❦","20250725T162518.000Z")
Add-ScrimRow $wsRoF $rof_T1 64 @("DRACO","KIT","MEEPLE","R-T","BONNIE","ALLI","Equipo 1","LOUD|FireCrow","LOUD|Edinho","LOUD|KaioDog","GO|Yichy❦","Tilo🍥","Bielz","20250725T161912.000Z")
Add-ScrimRow $wsRoF $rof_T1 65 @("DRACO","KIT","MEEPLE","R-T","BONNIE","ALLI","Equipo 1","LOUD|FireCrow","LOUD|Edinho","LOUD|KaioDog","GO|Yichy❦","Tilo🍥","Bielz","20250725T161652.000Z")
